$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A156:G156").Copy($ws.Range("A157:G157"))
$ws.Cells.Item(157,1).Value = 44433
$ws.Cells.Item(157,2).Value = 1482236
$ws.Cells.Item(157,3).Value = 1335021
$ws.Cells.Item(157,4).Value = 667299
$ws.Cells.Item(157,5).Value = 457579
$ws.Cells.Item(157,6).Value = 39750
$ws.Cells.Item(157,7).Value = 39485
$ws.Range("A157:G157").Copy($ws.Range("A158:G158"))
$ws.Cells.Item(158,1).Value = 44434
$ws.Cells.Item(158,2).Value = 1482236
$ws.Cells.Item(158,3).Value = 1349019
$ws.Cells.Item(158,4).Value = 667299
$ws.Cells.Item(158,5).Value = 466482
$ws.Cells.Item(158,6).Value = 39750
$ws.Cells.Item(158,7).Value = 39485
$ws.Range("A158:G158").Copy($ws.Range("A159:G159"))
$ws.Cells.Item(159,1).Value = 44435
$ws.Cells.Item(159,2).Value = 1482236
$ws.Cells.Item(159,3).Value = 1357601
$ws.Cells.Item(159,4).Value = 667299
$ws.Cells.Item(159,5).Value = 475016
$ws.Cells.Item(159,6).Value = 39750
$ws.Cells.Item(159,7).Value = 39485
$ws.Range("A159:G159").Copy($ws.Range("A160:G160"))
$ws.Cells.Item(160,1).Value = 44436
$ws.Cells.Item(160,2).Value = 1482236
$ws.Cells.Item(160,3).Value = 1400631
$ws.Cells.Item(160,4).Value = 667299
$ws.Cells.Item(160,5).Value = 478997
$ws.Cells.Item(160,6).Value = 39750
$ws.Cells.Item(160,7).Value = 39942
$ws.Range("A160:G160").Copy($ws.Range("A161:G161"))
$ws.Cells.Item(161,1).Value = 44438
$ws.Cells.Item(161,2).Value = 1482236
$ws.Cells.Item(161,3).Value = 1417132
$ws.Cells.Item(161,4).Value = 667299
$ws.Cells.Item(161,5).Value = 490232
$ws.Cells.Item(161,6).Value = 39750
$ws.Cells.Item(161,7).Value = 39952
$ws.Range("A161:G161").Copy($ws.Range("A162:G162"))
$ws.Cells.Item(162,1).Value = 44439
$ws.Cells.Item(162,2).Value = 1505856
$ws.Cells.Item(162,3).Value = 1420471
$ws.Cells.Item(162,4).Value = 800786
$ws.Cells.Item(162,5).Value = 497561
$ws.Cells.Item(162,6).Value = 39750
$ws.Cells.Item(162,7).Value = 39952
$ws.Range("A162:G162").Copy($ws.Range("A163:G163"))
$ws.Cells.Item(163,1).Value = 44440
$ws.Cells.Item(163,2).Value = 1505856
$ws.Cells.Item(163,3).Value = 1426524
$ws.Cells.Item(163,4).Value = 800786
$ws.Cells.Item(163,5).Value = 509292
$ws.Cells.Item(163,6).Value = 39750
$ws.Cells.Item(163,7).Value = 39952
$ws.Range("A163:G163").Copy($ws.Range("A164:G164"))
$ws.Cells.Item(164,1).Value = 44441
$ws.Cells.Item(164,2).Value = 1505856
$ws.Cells.Item(164,3).Value = 1430500
$ws.Cells.Item(164,4).Value = 800786
$ws.Cells.Item(164,5).Value = 526185
$ws.Cells.Item(164,6).Value = 39750
$ws.Cells.Item(164,7).Value = 39952
$ws.Range("A164:G164").Copy($ws.Range("A165:G165"))
$ws.Cells.Item(165,1).Value = 44442
$ws.Cells.Item(165,2).Value = 1505856
$ws.Cells.Item(165,3).Value = 1433995
$ws.Cells.Item(165,4).Value = 800786
$ws.Cells.Item(165,5).Value = 535772
$ws.Cells.Item(165,6).Value = 39750
$ws.Cells.Item(165,7).Value = 39952
$ws.Range("A165:G165").Copy($ws.Range("A166:G166"))
$ws.Cells.Item(166,1).Value = 44443
$ws.Cells.Item(166,2).Value = 1505856
$ws.Cells.Item(166,3).Value = 1434678
$ws.Cells.Item(166,4).Value = 800786
$ws.Cells.Item(166,5).Value = 538935
$ws.Cells.Item(166,6).Value = 39750
$ws.Cells.Item(166,7).Value = 39952
$ws.Range("A166:G166").Copy($ws.Range("A167:G167"))
$ws.Cells.Item(167,1).Value = 44445
$ws.Cells.Item(167,2).Value = 1505856
$ws.Cells.Item(167,3).Value = 1435077
$ws.Cells.Item(167,4).Value = 800786
$ws.Cells.Item(167,5).Value = 552857
$ws.Cells.Item(167,6).Value = 39750
$ws.Cells.Item(167,7).Value = 39952
$ws.Range("A167:G167").Copy($ws.Range("A168:G168"))
$ws.Cells.Item(168,1).Value = 44447
$ws.Cells.Item(168,2).Value = 1505856
$ws.Cells.Item(168,3).Value = 1436753
$ws.Cells.Item(168,4).Value = 800786
$ws.Cells.Item(168,5).Value = 565854
$ws.Cells.Item(168,6).Value = 39750
$ws.Cells.Item(168,7).Value = 39952
$ws.Range("A168:G168").Copy($ws.Range("A169:G169"))
$ws.Cells.Item(169,1).Value = 44448
$ws.Cells.Item(169,2).Value = 1527766
$ws.Cells.Item(169,3).Value = 1443028
$ws.Cells.Item(169,4).Value = 914366
$ws.Cells.Item(169,5).Value = 588034
$ws.Cells.Item(169,6).Value = 39750
$ws.Cells.Item(169,7).Value = 39952
$ws.Range("A169:G169").Copy($ws.Range("A170:G170"))
$ws.Cells.Item(170,1).Value = 44449
$ws.Cells.Item(170,2).Value = 1527766
$ws.Cells.Item(170,3).Value = 1451814
$ws.Cells.Item(170,4).Value = 914366
$ws.Cells.Item(170,5).Value = 602888
$ws.Cells.Item(170,6).Value = 39750
$ws.Cells.Item(170,7).Value = 39952
$ws.Range("A170:G170").Copy($ws.Range("A171:G171"))
$ws.Cells.Item(171,1).Value = 44450
$ws.Cells.Item(171,2).Value = 1527766
$ws.Cells.Item(171,3).Value = 1453771
$ws.Cells.Item(171,4).Value = 914366
$ws.Cells.Item(171,5).Value = 606539
$ws.Cells.Item(171,6).Value = 39750
$ws.Cells.Item(171,7).Value = 39952
$ws.Range("A171:G171").Copy($ws.Range("A172:G172"))
$ws.Cells.Item(172,1).Value = 44452
$ws.Cells.Item(172,2).Value = 1527766
$ws.Cells.Item(172,3).Value = 1462728
$ws.Cells.Item(172,4).Value = 914366
$ws.Cells.Item(172,5).Value = 634665
$ws.Cells.Item(172,6).Value = 39750
$ws.Cells.Item(172,7).Value = 39952

$ws.Range("D168").Select()
